$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.768.36'
$ws.Range('E2').Value = '  -2.22%  '
$ws.Range('D3').Value = '3.372.23'
$ws.Range('E3').Value = '  -4.06%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''555.55'
$ws.Range('E5').Value = '  -5.19%  '
$ws.Range('D6').Value = '''176.48'
$ws.Range('E6').Value = '  -1.53%  '
$ws.Range('E7').Value = '  -2.83%  '
$ws.Range('D8').Value = '3.363.97'
$ws.Range('E8').Value = '  -4.06%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = '''0.629'
$ws.Range('E10').Value = '  -1.79%  '
$ws.Range('D11').Value = '''0.163'
$ws.Range('E11').Value = '  -0.76%  '
$ws.Range('D12').Value = '''55.25'
$ws.Range('E12').Value = '  -0.98%  '
$ws.Range('D13').Value = '''0.0000273'
$ws.Range('E13').Value = '  -2.59%  '
$ws.Range('D14').Value = '''9.08'
$ws.Range('E14').Value = '  -2.49%  '
$ws.Range('D15').Value = '3.903.69'
$ws.Range('E15').Value = '  -4.23%  '
$ws.Range('D16').Value = '''18.38'
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = '''0.118'
$ws.Range('E17').Value = '  -2.67%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.364.20'
$ws.Range('E18').Value = '  -4.34%  '
$ws.Range('E19').Value = '  -2.09%  '
$ws.Range('D20').Value = '64.603.16'
$ws.Range('E20').Value = '  -2.46%  '
$ws.Range('D21').Value = '''0.982'
$ws.Range('E21').Value = '  -3.25%  '
$ws.Range('D22').Value = '''436.41'
$ws.Range('E22').Value = '  +4.95%  '
$ws.Range('D23').Value = '''4.99'
$ws.Range('E23').Value = '  +12.10%  '
$ws.Range('D24').Value = '''4.09'
$ws.Range('E24').Value = '  -5.50%  '
$ws.Range('D25').Value = '''84.49'
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('D26').Value = '''13.28'
$ws.Range('E26').Value = '  -2.25%  '
$ws.Range('D27').Value = '''10.84'
$ws.Range('E27').Value = '  -2.65%  '
$ws.Range('D28').Value = '''2.84'
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('D29').Value = '''8.79'
$ws.Range('E29').Value = '  -4.43%  '
$ws.Range('D30').Value = '''29.81'
$ws.Range('E30').Value = '  -2.09%  '
$ws.Range('D31').Value = '''6.64'
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('D32').Value = '''11.49'
$ws.Range('E32').Value = '  -2.78%  '
$ws.Range('D33').Value = '''582.04'
$ws.Range('E33').Value = '  -4.20%  '
$ws.Range('E34').Value = '  -3.09%  '
$ws.Range('D35').Value = '''58.69'
$ws.Range('E35').Value = '  -3.02%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  -7.27%  '
$ws.Range('D38').Value = '''3.54'
$ws.Range('E38').Value = '  -3.55%  '
$ws.Range('D39').Value = '''35.82'
$ws.Range('E39').Value = '  -3.17%  '
$ws.Range('D40').Value = '0.0₃0760'
$ws.Range('E40').Value = '  -5.45%  '
$ws.Range('D41').Value = '''0.369'
$ws.Range('E41').Value = '  -4.41%  '
$ws.Range('D42').Value = '3.117.03'
$ws.Range('E42').Value = '  -4.39%  '
$ws.Range('D43').Value = '''0.998'
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('E44').Value = '  -5.62%  '
$ws.Range('D45').Value = '''3.29'
$ws.Range('E45').Value = '  -1.82%  '
$ws.Range('D46').Value = '''0.0411'
$ws.Range('E46').Value = '  -2.79%  '
$ws.Range('E47').Value = '  -3.45%  '
$ws.Range('E48').Value = '  -2.75%  '
$ws.Range('E49').Value = '  -4.13%  '
$ws.Range('D50').Value = '''8.33'
$ws.Range('E50').Value = '  -3.97%  '
$ws.Range('D51').Value = '''134.99'
$ws.Range('E51').Value = '  -3.61%  '
